$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets hold the same event table and need the
# same "想去人数" (want-to-go count) updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 21
    $ws.Range("F4").Value = 45
    $ws.Range("F5").Value = 9
}
